$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 13 (shifts existing row 13 "Cristina Bormio Nunes" down to row 14)
$ws.Rows.Item(13).Insert()

# Populate the new row 13 with the additional professor
$ws.Range("B13").Value = "5840730 - Antonio Jefferson da Silva Machado"
$ws.Range("C13").Value = "5840730 - Antonio Jefferson da Silva Machado"

# Match styles used by the analogous rows (B column style index 2, C column style index 3)
$ws.Range("B13").Style = $ws.Range("B14").Style
$ws.Range("C13").Style = $ws.Range("C14").Style
